$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the contents of the two observation records currently
# stored in rows 3 and 4 (same underlying data, rows exchanged), while a
# handful of columns (C, D, I, J-O, T-V, W, X, Y, AD, AE, AF, AG, AT, AY, ...)
# hold identical values on both rows and are therefore left untouched.

# Capture the "before" values for every column that differs between the two
# rows, including the Starttid/Sluttid (Z/AB) vs public comment (AC) pair,
# which exists on only one of the two rows at a time.
$cols = @("A","B","E","F","G","H","P","Q","R","S","Z","AB","AC","AW","AX")

$row3vals = @{}
$row4vals = @{}
foreach ($col in $cols) {
    $row3vals[$col] = $ws.Range($col + "3").Value2
    $row4vals[$col] = $ws.Range($col + "4").Value2
}

foreach ($col in $cols) {
    $target3 = $ws.Range($col + "3")
    $target4 = $ws.Range($col + "4")

    if ($row4vals[$col] -eq $null) {
        $target3.ClearContents()
    } else {
        $target3.Value = $row4vals[$col]
    }

    if ($row3vals[$col] -eq $null) {
        $target4.ClearContents()
    } else {
        $target4.Value = $row3vals[$col]
    }
}
